# Add a new "Save" column (H) to the s_vals sheet, matching the
# formatting of the existing header cells and adding a numeric 0 value
# for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the
# neighboring "sum" header cell (G1) onto the new "Save" header cell (H1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data value for the new column in row 2.
$ws.Range("H2").Value = 0
